# Scheduled runner update: refresh market-price derived columns (H-N) across
# the ALC/ARM/BSM/CRP/CUL sheets with the latest Universalis price pull.
$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H2").Value = 72797.57000000001
$ws.Range("I2").Value = 101762.6
$ws.Range("J2").Value = 385
$ws.Range("K2").Value = 101762.6
$ws.Range("L2").Value = 385
$ws.Range("M2").Value = -101649.6
$ws.Range("N2").Value = -611

$ws.Range("H4").Value = 362.66666
$ws.Range("I4").Value = 362.66666
$ws.Range("J4").Value = 0
$ws.Range("K4").Value = 362.66666
$ws.Range("L4").Value = 0
$ws.Range("M4").Value = -248.66666
$ws.Range("N4").ClearContents()

$ws.Range("H15").Value = 306.4
$ws.Range("I15").Value = 306.4
$ws.Range("K15").Value = 919.1999999999999
$ws.Range("M15").Value = -750.1999999999999

$ws.Range("H51").Value = 3738.7083
$ws.Range("I51").Value = 4244.6665
$ws.Range("J51").Value = 3435.1333
$ws.Range("K51").Value = 4244.6665
$ws.Range("L51").Value = 3435.1333
$ws.Range("M51").Value = -3760.6665
$ws.Range("N51").Value = -4403.1333

$ws.Range("H55").Value = 338.55554
$ws.Range("I55").Value = 292.42856
$ws.Range("J55").Value = 500
$ws.Range("K55").Value = 292.42856
$ws.Range("L55").Value = 500
$ws.Range("M55").Value = -78.42856
$ws.Range("N55").Value = -928

$ws.Range("H64").Value = 4528.8
$ws.Range("I64").Value = 4441.5835
$ws.Range("J64").Value = 4609.3076
$ws.Range("K64").Value = 4441.5835
$ws.Range("L64").Value = 4609.3076
$ws.Range("M64").Value = -4193.5835
$ws.Range("N64").Value = -5105.3076

$ws.Range("H67").Value = 4528.8
$ws.Range("I67").Value = 4441.5835
$ws.Range("J67").Value = 4609.3076
$ws.Range("K67").Value = 4441.5835
$ws.Range("L67").Value = 4609.3076
$ws.Range("M67").Value = -3583.5835
$ws.Range("N67").Value = -6325.3076

$ws.Range("H116").Value = 3251.913
$ws.Range("J116").Value = 3115.1667
$ws.Range("L116").Value = 3115.1667
$ws.Range("N116").Value = -9999.1667

$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H2").Value = 2511.7778
$ws.Range("I2").Value = 2555.6667
$ws.Range("J2").Value = 2424
$ws.Range("K2").Value = 2555.6667
$ws.Range("L2").Value = 2424
$ws.Range("M2").Value = -2442.6667
$ws.Range("N2").Value = -2650

$ws.Range("H45").Value = 1199.2632
$ws.Range("I45").Value = 982.25
$ws.Range("K45").Value = 982.25
$ws.Range("M45").Value = -605.25

$ws.Range("H88").Value = 2209.0908
$ws.Range("I88").Value = 2063.158
$ws.Range("J88").Value = 2407.1428
$ws.Range("K88").Value = 2063.158
$ws.Range("L88").Value = 2407.1428
$ws.Range("M88").Value = -1657.158
$ws.Range("N88").Value = -3219.1428

$ws.Range("H91").Value = 2209.0908
$ws.Range("I91").Value = 2063.158
$ws.Range("J91").Value = 2407.1428
$ws.Range("K91").Value = 2063.158
$ws.Range("L91").Value = 2407.1428
$ws.Range("M91").Value = -659.1579999999999
$ws.Range("N91").Value = -5215.1428

$ws.Range("H116").Value = 2511.7778
$ws.Range("I116").Value = 2555.6667
$ws.Range("J116").Value = 2424
$ws.Range("K116").Value = 2555.6667
$ws.Range("L116").Value = 2424
$ws.Range("M116").Value = -261.6667000000002
$ws.Range("N116").Value = -7012

$ws.Range("H122").Value = 3058.182
$ws.Range("I122").Value = 2998.7368
$ws.Range("K122").Value = 8996.2104
$ws.Range("M122").Value = -6546.2104

$ws.Range("H132").Value = 27824.275
$ws.Range("I132").Value = 32908.242
$ws.Range("J132").Value = 3857
$ws.Range("K132").Value = 98724.726
$ws.Range("L132").Value = 11571
$ws.Range("M132").Value = -96194.726
$ws.Range("N132").Value = -16631

$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H3").Value = 2511.7778
$ws.Range("I3").Value = 2555.6667
$ws.Range("J3").Value = 2424
$ws.Range("K3").Value = 2555.6667
$ws.Range("L3").Value = 2424
$ws.Range("M3").Value = -2441.6667
$ws.Range("N3").Value = -2652

$ws.Range("H86").Value = 3331.125
$ws.Range("I86").Value = 3992.1904
$ws.Range("J86").Value = 2069.0908
$ws.Range("K86").Value = 3992.1904
$ws.Range("L86").Value = 2069.0908
$ws.Range("M86").Value = -2869.1904
$ws.Range("N86").Value = -4315.0908

$ws.Range("H89").Value = 3331.125
$ws.Range("I89").Value = 3992.1904
$ws.Range("J89").Value = 2069.0908
$ws.Range("K89").Value = 19960.952
$ws.Range("L89").Value = 10345.454
$ws.Range("M89").Value = -14344.952
$ws.Range("N89").Value = -21577.454

$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H31").Value = 2364.3655
$ws.Range("I31").Value = 1618.775
$ws.Range("J31").Value = 4849.6665
$ws.Range("K31").Value = 1618.775
$ws.Range("L31").Value = 4849.6665
$ws.Range("M31").Value = -1323.775
$ws.Range("N31").Value = -5439.6665

$ws.Range("H34").Value = 2364.3655
$ws.Range("I34").Value = 1618.775
$ws.Range("J34").Value = 4849.6665
$ws.Range("K34").Value = 1618.775
$ws.Range("L34").Value = 4849.6665
$ws.Range("M34").Value = -1416.775
$ws.Range("N34").Value = -5253.6665

$ws.Range("H58").Value = 4615.147
$ws.Range("I58").Value = 5855.75
$ws.Range("J58").Value = 2842.8572
$ws.Range("K58").Value = 5855.75
$ws.Range("L58").Value = 2842.8572
$ws.Range("M58").Value = -5652.75
$ws.Range("N58").Value = -3248.8572

$ws.Range("H99").Value = 54681.316
$ws.Range("I99").Value = 168289.67
$ws.Range("J99").Value = 2246.6924
$ws.Range("K99").Value = 168289.67
$ws.Range("L99").Value = 2246.6924
$ws.Range("M99").Value = -166791.67
$ws.Range("N99").Value = -5242.6924

$ws.Range("H126").Value = 54681.316
$ws.Range("I126").Value = 168289.67
$ws.Range("J126").Value = 2246.6924
$ws.Range("K126").Value = 504869.01
$ws.Range("L126").Value = 6740.0772
$ws.Range("M126").Value = -502399.01
$ws.Range("N126").Value = -11680.0772

$ws.Range("H135").Value = 56743.25
$ws.Range("J135").Value = 39332.668
$ws.Range("L135").Value = 39332.668
$ws.Range("N135").Value = -49472.668

$ws.Range("H136").Value = 4615.147
$ws.Range("I136").Value = 5855.75
$ws.Range("J136").Value = 2842.8572
$ws.Range("K136").Value = 17567.25
$ws.Range("L136").Value = 8528.571599999999
$ws.Range("M136").Value = -15017.25
$ws.Range("N136").Value = -13628.5716

$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H113").Value = 13158443
$ws.Range("I113").Value = 20000560
$ws.Range("J113").Value = 527.38464
$ws.Range("K113").Value = 60001680
$ws.Range("L113").Value = 1582.15392
$ws.Range("M113").Value = -59999510
$ws.Range("N113").Value = -5922.15392

$ws.Range("H132").Value = 5379.5557
$ws.Range("I132").Value = 2358
$ws.Range("J132").Value = 10127.714
$ws.Range("K132").Value = 21222
$ws.Range("L132").Value = 91149.42600000001
$ws.Range("M132").Value = -18692
$ws.Range("N132").Value = -96209.42600000001

